$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new row at position 8, pushing existing rows 8-10 down to 9-11
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the "Betty" / database creation entry
$ws.Cells.Item(8, 1).Value = "Creación de la base de datos"
$ws.Cells.Item(8, 2).Value = "Betty"
$ws.Cells.Item(8, 4).Value = "10:45pm"
$ws.Cells.Item(8, 5).Value = "11:32pm"
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1

$ws.Cells.Item(8, 3).Value = [DateTime]"2015-04-25"

$ws.Range("G10").Select()
